$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# "Experimental" row (row 7): set the Value column to the literal text "true".
# A bare Value = "true" gets auto-coerced to a real Boolean by Excel, so force
# literal text with a leading apostrophe, then restore the original (non
# quote-prefixed) cell format by copying it from the row above.
$c = $ws.Range("B7")
$c.Value = "'true"
$ws.Range("B6").Copy()
$c.PasteSpecial(-4122)
$excel.CutCopyMode = $false

# "Date" row (row 8): update the timestamp value
$ws.Range("B8").Value = "2025-07-21T12:46:15+00:00"
